$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 655, shifting existing rows (655:781) down to (656:782)
$ws.Rows.Item(655).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 655 with the new data record
$ws.Cells.Item(655, 1).Value = 10
$ws.Cells.Item(655, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(655, 3).Value = "La Araucanía"
$ws.Cells.Item(655, 4).Value = 45209
$ws.Cells.Item(655, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(655, 5).Value = 9
$ws.Cells.Item(655, 6).Value = 100112023
$ws.Cells.Item(655, 7).Value = "Brócoli"
$ws.Cells.Item(655, 8).Value = "Sin especificar"
$ws.Cells.Item(655, 9).Value = "Primera"
$ws.Cells.Item(655, 10).Value = 500
$ws.Cells.Item(655, 11).Value = 1200
$ws.Cells.Item(655, 12).Value = 1200
$ws.Cells.Item(655, 13).Value = 1200
$ws.Cells.Item(655, 14).Value = "$/unidad"
$ws.Cells.Item(655, 15).Value = "Región Metropolitana"
$ws.Cells.Item(655, 16).Value = 1200
$ws.Cells.Item(655, 17).Value = 1
$ws.Cells.Item(655, 18).Value = "Hortaliza"
